$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal string (e.g. "145.20") must be forced
# to Text so Excel does not silently coerce them to a Number (losing trailing
# zeros / introducing float rounding). We flip NumberFormat to Text, write the
# value, then reset the style back to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = "62.609.92"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.439.59"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "2.437.04"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.06%  "
$ws.Range("E15").Value = "  +4.89%  "
$ws.Range("D16").Value = "2.883.23"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "62.520.98"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "2.444.04"
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").Value = "0.0₆0803"
$ws.Range("E21").Value = "  +180.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "326.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  +10.37%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "626.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.33%  "
$ws.Range("E28").Value = "  +12.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("E30").Value = "  +4.32%  "
$ws.Range("D31").Value = "2.565.16"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("E33").Value = "  +5.83%  "
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.372"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "151.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.39%  "
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +28.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.601"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.95%  "
